$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A, B, and the general/default-ish C.. columns) ---
$ws.Columns("A").ColumnWidth = 41.25
$ws.Columns("B").ColumnWidth = 8.75
$ws.Range("C1:I1").EntireColumn.ColumnWidth = 41.25

# --- Content: title cell + six explanatory rows (merged across columns) ---
$ws.Range("A46").Value = ".gitignore对.idea文件忽略没效果"

$ws.Range("C46:I46").Merge()
$ws.Range("C46").Value = "在.idea文件目录运行git.bash"

$ws.Range("C47:I47").Merge()
$ws.Range("C47").Value = "git rm -r --cached .idea/*"

$ws.Range("C48:I48").Merge()
$ws.Range("C48").Value = "上句代码会删除掉缓存中的.idea下所有文件，实际目录没有"

$ws.Range("C49:I49").Merge()
$ws.Range("C49").Value = "按改动文件正常提交到origin/master"

$ws.Range("C50:I50").Merge()
$ws.Range("C50").Value = "在.gitignore文件中添加.idea并提交到origin/master即可"

$ws.Range("C51:G51").Merge()
$ws.Range("C51").Value = "【原因就是git已经关联追踪了这些文件，再次设置ignore时无效。需要执行一次 git rm -r --cached 目录名称/文件名称】"

# --- Fonts: 16pt 微软雅黑 (theme text colour) everywhere in this block ---
$ws.Range("A46").Font.Name = "微软雅黑"
$ws.Range("A46").Font.Size = 16
$ws.Range("A46").Font.ThemeColor = 1

$ws.Range("C46:I50").Font.Name = "微软雅黑"
$ws.Range("C46:I50").Font.Size = 16
$ws.Range("C46:I50").Font.ThemeColor = 1
$ws.Range("C46:I50").HorizontalAlignment = -4131

$ws.Range("C51:G51").Font.Name = "微软雅黑"
$ws.Range("C51:G51").Font.Size = 16
$ws.Range("C51:G51").Font.ThemeColor = 1
$ws.Range("C51:G51").HorizontalAlignment = -4131

# --- Sheet view / selection ---
$ws.Range("C51:G51").Select()

# --- Page setup ---
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
